# Apply edits described by the commit: "added tweet with photo & tweet with image count"
#  1) Update a handful of 'languages' row-6 cells whose underlying python-side
#     langdetect sets were recomputed.
#  2) Append two new summary rows: hasImages (row 8) and hasVideos (row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row 6 ("languages") cell corrections -------------------------------
$ws.Range("D6").Value = "['ar' 'ca' 'es' 'so' 'id']"
$ws.Range("F6").Value = "['ar' 'undetected' 'en']"
$ws.Range("K6").Value = "['en' 'ar' 'ja' 'de' 'id' 'pt' 'fr']"
$ws.Range("O6").Value = "['pl' 'en' 'ar' 'undetected' 'id']"
$ws.Range("Q6").Value = "['ar' 'es' 'fr' 'en' 'so' 'undetected']"
$ws.Range("R6").Value = "['ar' 'en' 'fr' 'undetected' 'nl' 'tr' 'id' 'ca' 'es']"
$ws.Range("U6").Value = "['ar' 'es' 'en' 'undetected' 'pt' 'hu' 'fa' 'ca' 'id' 'et']"
$ws.Range("Y6").Value = "['ca' 'ar' 'es' 'undetected' 'en']"

# --- 2) New row 8: hasImages -------------------------------------------------
$ws.Range("A8").Value = "hasImages"

$hasImages = @(20, 4, 7, 42, 38, 12, 66, 1305, 14, 108, 44, 30, 13, 55, 15, 251, 168, 43, 55, 187, 443, 169, 22, 12)
for ($i = 0; $i -lt $hasImages.Length; $i++) {
    $ws.Cells.Item(8, $i + 2).Value = $hasImages[$i]
}

# --- 3) New row 9: hasVideos -------------------------------------------------
$ws.Range("A9").Value = "hasVideos"

$hasVideos = @(12, 0, 6, 5, 8, 1, 14, 140, 13, 11, 9, 13, 0, 0, 5, 62, 39, 12, 2, 91, 64, 11, 0, 6)
for ($i = 0; $i -lt $hasVideos.Length; $i++) {
    $ws.Cells.Item(9, $i + 2).Value = $hasVideos[$i]
}

# --- 4) Match the bold/centered/bordered label style used by column A ------
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
